$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.972.75"
$ws.Range("E2").Value = "  +4.48%  "
$ws.Range("D3").Value = "1.780.05"
$ws.Range("E3").Value = "  +2.89%  "
$ws.Range("D4").Value = "'0.9981"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'242.91"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "'0.9986"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.4880"
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("D8").Value = "'0.2658"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("D9").Value = "'0.06240"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "1.780.63"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("D11").Value = "'16.33"
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("D12").Value = "'0.06987"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'0.6169"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "'4.591"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").Value = "'79.47"
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "27.920.50"
$ws.Range("E16").Value = "  +5.14%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'0.9980"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'0.000007205"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").Value = "'11.75"
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("D21").Value = "2.008.22"
$ws.Range("E21").Value = "  +2.80%  "
$ws.Range("D22").Value = "'4.577"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").Value = "'8.660"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").Value = "'5.183"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "'141.87"
$ws.Range("D26").Value = "'15.57"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").Value = "'1.872"
$ws.Range("E27").Value = "  +6.53%  "
$ws.Range("D28").Value = "'109.46"
$ws.Range("E28").Value = "  +3.09%  "
$ws.Range("D29").Value = "'1.402"
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("D30").Value = "'4.072"
$ws.Range("E30").Value = "  +3.04%  "
$ws.Range("D31").Value = "'0.08265"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").Value = "'3.770"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("D33").Value = "'0.04718"
$ws.Range("E33").Value = "  +4.10%  "
$ws.Range("D34").Value = "'1.056"
$ws.Range("E34").Value = "  +5.26%  "
$ws.Range("D35").Value = "'2.602"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "'0.6336"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "'0.9387"
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("D38").Value = "'2.584"
$ws.Range("E38").Value = "  +6.81%  "
$ws.Range("D39").Value = "'2.058"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").Value = "'5.838"
$ws.Range("D41").Value = "'0.01533"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'100.06"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "'0.3928"
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").Value = "'7.148"
$ws.Range("E45").Value = "  +2.52%  "
$ws.Range("D46").Value = "'0.1191"
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("D47").Value = "'0.05405"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "'7.910"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("D49").Value = "'30.36"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").Value = "'1.275"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("D51").Value = "'52.33"
$ws.Range("E51").Value = "  +1.14%  "
